$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: force Word to keep two adjacent same-formatted runs distinct
# (rather than merging them into a single run) by toggling Bold on then
# back off across the given range. The net formatting effect is a no-op
# but the run boundary survives the save/normalize pass.
# ---------------------------------------------------------------------

# 1. "Add a thread group" -> append " and give some users to it" as its
#    own run (paragraph 3).
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$r3.End = $r3.End - 1
$insStart = $r3.End
$r3.InsertAfter(" and give some users to it")
$newRun = $d.Range($insStart, $insStart + (" and give some users to it").Length)
$newRun.Bold = 1
$newRun.Bold = 0

# 2. Picture paragraph right after it gets NoProofing (w:noProof).
$d.Paragraphs(4).Range.NoProofing = $true

# 3. "Add a throughput " + proofErr + "controllers" + proofErr ->
#    single run "Add a throughput controllers" with the proofErr marks
#    gone. Delete the whole paragraph (this drops the proofErr nodes
#    too) and recreate a fresh one in its place.
$p5 = $d.Paragraphs(5)
$full5 = $d.Range($p5.Range.Start, $p5.Range.End)
$full5.Delete()
$p5recreated = $d.Paragraphs(5)
$p5recreated.Range.InsertParagraphBefore()
$p5target = $d.Paragraphs(5)
$r5 = $p5target.Range
$r5.End = $r5.End - 1
$r5.InsertAfter("Add a throughput controllers")

# 4. Picture paragraph right after it gets NoProofing.
$d.Paragraphs(6).Range.NoProofing = $true

# 5. "Add http requests under that controller" unchanged (paragraph 7).

# 6. Picture paragraph right after it gets NoProofing.
$d.Paragraphs(8).Range.NoProofing = $true

# 7. "Throughput controller is used to define number of times the test
#    can run by giving it a percentage value" -> split into three runs,
#    replacing "percentage" with "number of users".
$found9 = $d.Content.Find.Execute("percentage", $true, $false, $false, $false, $false, $true, 1, $false, "number of users", 2)
$p9 = $d.Paragraphs(9)
$full9text = $p9.Range.Text
$offset9 = $full9text.IndexOf("number of users")
$rStart9 = $p9.Range.Start + $offset9
$rEnd9 = $rStart9 + ("number of users").Length
$sub9 = $d.Range($rStart9, $rEnd9)
$sub9.Bold = 1
$sub9.Bold = 0

# 8. "Adding a listener to see the output" unchanged (paragraph 10).
# 9. "Aggregate report" unchanged (paragraph 11).

# 10. Picture paragraph after "Aggregate report" gets NoProofing.
$d.Paragraphs(12).Range.NoProofing = $true

# 11. Empty paragraph 13 unchanged.

# 12. Picture paragraph (14) gets NoProofing.
$d.Paragraphs(14).Range.NoProofing = $true

# 13. "We have added 3 throughput controllers ..." unchanged (paragraph 15).

# 14. "So that the requests are executed according to the percentage
#     value" -> becomes "Total users are 100, we have given the users
#     to request in the ratio 50:30:20", and a NEW paragraph is added
#     after it with the old sentence (split into 3 runs, "percentage"
#     replaced by "number of users").
$p16 = $d.Paragraphs(16)
$r16 = $p16.Range
$r16.End = $r16.End - 1
$r16.Text = "__PLACEHOLDER_FOR_REWRITE__"
$r16b = $d.Paragraphs(16).Range
$r16b.End = $r16b.End - 1
$r16b.Text = "Total users are 100, we have given the users to request in the ratio 50:30:20"

$p16c = $d.Paragraphs(16)
$p16c.Range.InsertParagraphAfter()

$p17 = $d.Paragraphs(17)
$r17 = $p17.Range
$r17.End = $r17.End - 1
$r17.InsertAfter("So that the requests are executed according to the number of users value")

$full17text = $d.Paragraphs(17).Range.Text
$offset17 = $full17text.IndexOf("number of users")
$pStart17 = $d.Paragraphs(17).Range.Start
$rStart17 = $pStart17 + $offset17
$rEnd17 = $rStart17 + ("number of users").Length
$sub17 = $d.Range($rStart17, $rEnd17)
$sub17.Bold = 1
$sub17.Bold = 0

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
